$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 315-316; this pushes the existing rows
# 315-324 down to 317-326 (data itself stays exactly as-is, just moves).
$ws.Range("A315:A316").EntireRow.Insert()

# Populate the newly inserted row 315 with its data.
$ws.Cells.Item(315, 1).Value = 11
$ws.Cells.Item(315, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(315, 3).Value = "Bíobío"
$ws.Cells.Item(315, 4).Value = 44706
$ws.Cells.Item(315, 5).Value = 8
$ws.Cells.Item(315, 6).Value = 100112006
$ws.Cells.Item(315, 7).Value = "Repollo"
$ws.Cells.Item(315, 8).Value = "Crespo record"
$ws.Cells.Item(315, 9).Value = "Primera"
$ws.Cells.Item(315, 10).Value = 2000
$ws.Cells.Item(315, 11).Value = 900
$ws.Cells.Item(315, 12).Value = 1000
$ws.Cells.Item(315, 13).Value = 950
$ws.Cells.Item(315, 14).Value = "$/unidad"
$ws.Cells.Item(315, 15).Value = "Región Metropolitana"
$ws.Cells.Item(315, 16).Value = 950
$ws.Cells.Item(315, 17).Value = 1
$ws.Cells.Item(315, 18).Value = "Hortaliza"

# Populate the newly inserted row 316 with its data.
$ws.Cells.Item(316, 1).Value = 11
$ws.Cells.Item(316, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(316, 3).Value = "Bíobío"
$ws.Cells.Item(316, 4).Value = 44706
$ws.Cells.Item(316, 5).Value = 8
$ws.Cells.Item(316, 6).Value = 100112006
$ws.Cells.Item(316, 7).Value = "Repollo"
$ws.Cells.Item(316, 8).Value = "Crespo record"
$ws.Cells.Item(316, 9).Value = "Segunda"
$ws.Cells.Item(316, 10).Value = 1000
$ws.Cells.Item(316, 11).Value = 800
$ws.Cells.Item(316, 12).Value = 800
$ws.Cells.Item(316, 13).Value = 800
$ws.Cells.Item(316, 14).Value = "$/unidad"
$ws.Cells.Item(316, 15).Value = "Región Metropolitana"
$ws.Cells.Item(316, 16).Value = 800
$ws.Cells.Item(316, 17).Value = 1
$ws.Cells.Item(316, 18).Value = "Hortaliza"
